$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.363.68'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '2.315.37'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.07'
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.33'
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').Value = '2.337.42'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -1.54%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.32'
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').Value = '2.731.65'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.40'
$ws.Range('E15').Value = '  -4.23%  '
$ws.Range('D16').Value = '57.204.80'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '2.324.88'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '335.24'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.41'
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.80'
$ws.Range('E22').Value = '  +1.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.78'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.69'
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.76'
$ws.Range('E29').Value = '  +3.67%  '
$ws.Range('E30').Value = '  +0.80%  '
$ws.Range('D31').Value = '0.0₃0724'
$ws.Range('E31').Value = '  -3.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.10'
$ws.Range('E32').Value = '  -2.83%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('E36').Value = '  -3.69%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.98'
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.15'
$ws.Range('E39').Value = '  +1.00%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.57'
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.76'
$ws.Range('E41').Value = '  +8.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '148.05'
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.374'
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('E44').Value = '  -1.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '282.68'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('E48').Value = '  +3.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.557'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('E50').Value = '  +5.09%  '
$ws.Range('E51').Value = '  -1.57%  '
